$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 500
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1308

$ws.Range("H129").Value = 925.8472
$ws.Range("J129").Value = 955.4706
$ws.Range("L129").Value = 2866.4118
$ws.Range("N129").Value = -12866.4118

$ws.Range("H138").Value = 2974.5
$ws.Range("I138").Value = 2234.0715
$ws.Range("J138").Value = 3153.224
$ws.Range("K138").Value = 6702.2145
$ws.Range("L138").Value = 9459.672
$ws.Range("M138").Value = -1562.2145
$ws.Range("N138").Value = -19739.672

$ws.Range("H141").Value = 2041.9706
$ws.Range("I141").Value = 1764.0667
$ws.Range("J141").Value = 4126.25
$ws.Range("K141").Value = 5292.2001
$ws.Range("L141").Value = 12378.75
$ws.Range("M141").Value = -112.2001
$ws.Range("N141").Value = -22738.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1202.7273

$ws.Range("H5").Value = 189.22223
$ws.Range("I5").Value = 120.25
$ws.Range("J5").Value = 244.4
$ws.Range("K5").Value = 120.25
$ws.Range("L5").Value = 244.4
$ws.Range("M5").Value = -8.25
$ws.Range("N5").Value = -468.4

$ws.Range("H61").Value = 2333.5293
$ws.Range("I61").Value = 1628.7142
$ws.Range("J61").Value = 2826.9
$ws.Range("K61").Value = 1628.7142
$ws.Range("L61").Value = 2826.9
$ws.Range("M61").Value = -1416.7142
$ws.Range("N61").Value = -3250.9

$ws.Range("H116").Value = 1202.7273

$ws.Range("H122").Value = 3140.4
$ws.Range("I122").Value = 1822.6666
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 5467.9998
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -3017.9998
$ws.Range("N122").Value = -49900

$ws.Range("H132").Value = 2586.5518
$ws.Range("I132").Value = 1195.6
$ws.Range("J132").Value = 5677.5557
$ws.Range("K132").Value = 3586.8
$ws.Range("L132").Value = 17032.6671
$ws.Range("M132").Value = -1056.8
$ws.Range("N132").Value = -22092.6671

$ws.Range("H136").Value = 2333.5293
$ws.Range("I136").Value = 1628.7142
$ws.Range("J136").Value = 2826.9
$ws.Range("K136").Value = 4886.142599999999
$ws.Range("L136").Value = 8480.700000000001
$ws.Range("M136").Value = -2336.142599999999
$ws.Range("N136").Value = -13580.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1202.7273

$ws.Range("H4").Value = 189.22223
$ws.Range("I4").Value = 120.25
$ws.Range("J4").Value = 244.4
$ws.Range("K4").Value = 120.25
$ws.Range("L4").Value = 244.4
$ws.Range("M4").Value = -5.25
$ws.Range("N4").Value = -474.4

$ws.Range("H94").Value = 2089.7646
$ws.Range("I94").Value = 2135.0667
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 2135.0667
$ws.Range("L94").Value = 1750
$ws.Range("M94").Value = -1684.0667
$ws.Range("N94").Value = -2652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 33065.6
$ws.Range("I64").Value = 20246
$ws.Range("J64").Value = 36270.5
$ws.Range("K64").Value = 20246
$ws.Range("L64").Value = 36270.5
$ws.Range("M64").Value = -19998
$ws.Range("N64").Value = -36766.5

$ws.Range("H67").Value = 33065.6
$ws.Range("I67").Value = 20246
$ws.Range("J67").Value = 36270.5
$ws.Range("K67").Value = 20246
$ws.Range("L67").Value = 36270.5
$ws.Range("M67").Value = -19388
$ws.Range("N67").Value = -37986.5

$ws.Range("H99").Value = 6158
$ws.Range("I99").Value = 3403.6667
$ws.Range("J99").Value = 11666.667
$ws.Range("K99").Value = 3403.6667
$ws.Range("L99").Value = 11666.667
$ws.Range("M99").Value = -1905.6667
$ws.Range("N99").Value = -14662.667

$ws.Range("H122").Value = 2531.2942
$ws.Range("I122").Value = 1752
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 5256
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -2806
$ws.Range("N122").Value = -49900

$ws.Range("H126").Value = 6158
$ws.Range("I126").Value = 3403.6667
$ws.Range("J126").Value = 11666.667
$ws.Range("K126").Value = 10211.0001
$ws.Range("L126").Value = 35000.001
$ws.Range("M126").Value = -7741.000100000001
$ws.Range("N126").Value = -39940.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 4225
$ws.Range("J48").Value = 5950
$ws.Range("L48").Value = 17850
$ws.Range("N48").Value = -18350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4530.3076
$ws.Range("J122").Value = 15000
$ws.Range("L122").Value = 45000
$ws.Range("N122").Value = -49900

$ws.Range("H132").Value = 5494.9165
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 5085.364
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 15256.092
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -20316.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3150
$ws.Range("J46").Value = 2687.5
$ws.Range("L46").Value = 2687.5
$ws.Range("N46").Value = -3063.5

$ws.Range("H122").Value = 3389.825
$ws.Range("I122").Value = 3097.1843
$ws.Range("J122").Value = 8950
$ws.Range("K122").Value = 9291.552899999999
$ws.Range("L122").Value = 26850
$ws.Range("M122").Value = -6841.552899999999
$ws.Range("N122").Value = -31750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 938.1111
$ws.Range("I107").Value = 836.9231
$ws.Range("J107").Value = 1201.2
$ws.Range("K107").Value = 2510.7693
$ws.Range("L107").Value = 3603.6
$ws.Range("M107").Value = -590.7692999999999
$ws.Range("N107").Value = -7443.6

$ws.Range("H113").Value = 20401.2
$ws.Range("I113").Value = 25376.5
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 76129.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -73959.5
$ws.Range("N113").Value = -5840

$ws.Range("H122").Value = 5640.8
$ws.Range("I122").Value = 4600.8887
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 13802.6661
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -11352.6661
$ws.Range("N122").Value = -49900

$ws.Range("H132").Value = 4904672
$ws.Range("I132").Value = 3003.15
$ws.Range("J132").Value = 11907056
$ws.Range("K132").Value = 9009.450000000001
$ws.Range("L132").Value = 35721168
$ws.Range("M132").Value = -6479.450000000001
$ws.Range("N132").Value = -35726228
